$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.817.43'
$ws.Range('E2').Value = '  +0.05%  '
$ws.Range('D3').Value = '2.664.13'
$ws.Range('E3').Value = '  -0.61%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '597.52'
$ws.Range('D5').ClearFormats()
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '157.98'
$ws.Range('D6').ClearFormats()
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.652'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +4.31%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  -2.73%  '
$ws.Range('E10').Value = '  +0.35%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.86'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.41%  '
$ws.Range('E12').Value = '  +1.49%  '
$ws.Range('E13').Value = '  -1.52%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000195'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.28%  '
$ws.Range('D15').Value = '3.142.33'
$ws.Range('D16').Value = '65.663.88'
$ws.Range('E16').Value = '  +0.04%  '
$ws.Range('D17').Value = '2.656.47'
$ws.Range('E17').Value = '  -1.00%  '
$ws.Range('E18').Value = '  -2.00%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.81'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.18%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '351.03'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.46%  '
$ws.Range('E21').Value = '  -1.55%  '
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.84'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +11.80%  '
$ws.Range('E25').Value = '  +0.13%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.73'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.48%  '
$ws.Range('E27').Value = '  +1.18%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '566.23'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +6.41%  '
$ws.Range('E29').Value = '  +1.46%  '
$ws.Range('E30').Value = '  -2.99%  '
$ws.Range('E31').Value = '  -0.21%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.14'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.21%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.83'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +4.06%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.68'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.84%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.59'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.41%  '
$ws.Range('E36').Value = '  -0.69%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '20.63'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.20%  '
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.96'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.37%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '154.57'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -2.40%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '161.47'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.15%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '4.10'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.37%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0621'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.40%  '
$ws.Range('E44').Value = '  -1.07%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '23.05'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.40%  '
$ws.Range('E46').Value = '  +0.12%  '
$ws.Range('E47').Value = '  -0.21%  '
$ws.Range('E48').Value = '  +1.38%  '
$ws.Range('E49').Value = '  -1.69%  '
$ws.Range('E50').Value = '  -6.10%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.817'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.14%  '
